# Add a new row (71) of data to each of the four worksheets, mirroring the
# structure of the existing rows (time / length / id / actual-length /
# checksum plus their decimal counterparts).

$wb = $excel.ActiveWorkbook

$rowData = @{
    "DE_LFT_#1" = @{
        A = 45857.43685185185
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x48"
        E = "0x14"
        F = 380
        G = "7.598631275147109e+23"
        H = 328
        I = 14
    }
    "DE_LFT_#2" = @{
        A = 45857.43685185185
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x48"
        E = "0xe"
        F = 380
        G = "5.68432987514711e+23"
        H = 328
        I = 14
    }
    "DE_PLT_#1" = @{
        A = 45857.43685185185
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x7A"
        E = "0x7"
        F = 130
        G = "5.68631262647114e+23"
        H = 122
        I = 7
    }
    "DE_PLT_#2" = @{
        A = 45857.43685185185
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x7A"
        E = "0x3"
        F = 130
        G = "9.85046333984776e+23"
        H = 122
        I = 3
    }
}

foreach ($sheetName in $rowData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = 71
    $data = $rowData[$sheetName]

    # Column A keeps the same date/time number format used by the existing
    # rows in this column (numFmt "YYYY-MM-DD HH:MM:SS").
    $ws.Range("A$newRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("A$newRow").Value = $data.A

    # Columns B-E are text (hex byte listings), stored as strings.
    $ws.Range("B$newRow").Value = $data.B
    $ws.Range("C$newRow").Value = $data.C
    $ws.Range("D$newRow").Value = $data.D
    $ws.Range("E$newRow").Value = $data.E

    # Columns F-I are plain numeric decimal values.
    $ws.Range("F$newRow").Value = $data.F
    $ws.Range("G$newRow").Value = [double]$data.G
    $ws.Range("H$newRow").Value = $data.H
    $ws.Range("I$newRow").Value = $data.I
}
